$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (gpt-4o) updates
$ws.Range("B2").Value = 4
$ws.Range("H2").Value = 1
$ws.Range("N2").Value = 3.55
$ws.Range("O2").Value = "The report demonstrates strong evidence extraction quality with accurate citations and full sentence quotations, though it lacks some depth in coverage of representativeness dimensions, particularly in demographic and geographic specifics. The structure is clear and well-organized, aiding readability. The relevance and faithfulness of the evidence are maintained, with no unsupported assumptions. However, the identification of missing disclosures could be more detailed, particularly regarding specific demographic and geographic targets. The audit usefulness is reasonable, but the report could benefit from more explicit traceability and verifiability of claims. Overall, the report is strong but has clear areas for improvement, especially in detailing missing disclosures and enhancing audit traceability."

# Row 3 (ollama_mistral) updates
$ws.Range("F3").Value = 3
$ws.Range("L3").Value = 0.3
$ws.Range("N3").Value = 3.15
$ws.Range("O3").Value = "The evidence report is reasonably structured with clear sections and consistent formatting, earning a score of 4 in structure and formatting. However, the evidence extraction quality is only acceptable, as the report includes some fragmented snippets and lacks full sentence quotations, resulting in a score of 3. The coverage of representativeness dimensions is addressed but lacks depth and specificity, particularly in geographic and language distribution, leading to a score of 3. Relevance and faithfulness are maintained, but the report occasionally lacks direct grounding in the source documents, scoring a 3. Missing disclosures are identified, but not comprehensively, also scoring a 3. Overall, the report is useful for auditing but could be more detailed and precise, resulting in a score of 3 for audit usefulness."
